$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2-9 (columns A-T).
# Column layout:
# A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E Ligand-expressing cells, F Ligand detection rate,
# G Ligand average expression value, H Ligand total expression value,
# I Ligand derived specificity of average expression value,
# J Ligand derived specificity of total expression value,
# K Receptor-expressing cells, L Receptor detection rate,
# M Receptor average expression value, N Receptor total expression value,
# O Receptor derived specificity of average expression value,
# P Receptor derived specificity of total expression value,
# Q Edge average expression weight, R Edge total expression weight,
# S Edge average expression derived specificity, T Edge total expression derived specificity

$rows = @(
    @("FAPs","Lgi3","Adam22","ECs",3,1,1.535771,4.607313,0.8955152254457543,0.8955152254457543,3,1,3.078423333333333,9.23527,0.2543225865010234,0.2543225865010234,4.727753281056668,42.54977952951,0.2277497483864113,0.2277497483864113),
    @("FAPs","Lgi3","Adam22","FAPs",3,1,1.535771,4.607313,0.8955152254457543,0.8955152254457543,3,1,4.578777,13.736331,0.3782736432128339,0.3782736432128339,7.031952932067001,63.28757638860301,0.3387498068819278,0.3387498068819278),
    @("FAPs","Lgi3","Adam22","M2",3,1,1.535771,4.607313,0.8955152254457543,0.8955152254457543,3,1,0.8728403333333333,2.618521,0.07210931932983509,0.07210931932983509,1.340482871563667,12.064345844073,0.06457499335639717,0.06457499335639717),
    @("FAPs","Lgi3","Adam22","sCs",3,1,1.535771,4.607313,0.8955152254457543,0.8955152254457543,3,1,3.574363333333334,10.72309,0.2952944509563076,0.2952944509563076,5.489403550796668,49.40463195717001,0.2644406768210181,0.2644406768210181),
    @("sCs","Lgi3","Adam22","M2",1,0.3333333333333333,0.179187,0.537561,0.1044847745542456,0.1044847745542456,3,1,3.078423333333333,9.23527,0.2543225865010234,0.2543225865010234,0.5516134418299999,4.964520976469999,0.02657283811461206,0.02657283811461206),
    @("sCs","Lgi3","Adam22","FAPs",1,0.3333333333333333,0.179187,0.537561,0.1044847745542456,0.1044847745542456,3,1,4.578777,13.736331,0.3782736432128339,0.3782736432128339,0.8204573142989998,7.384115828690999,0.03952383633090609,0.03952383633090609),
    @("sCs","Lgi3","Adam22","M2",1,0.3333333333333333,0.179187,0.537561,0.1044847745542456,0.1044847745542456,3,1,0.8728403333333333,2.618521,0.07210931932983509,0.07210931932983509,0.156401640809,1.407614767281,0.007534325973437926,0.007534325973437926),
    @("sCs","Lgi3","Adam22","sCs",1,0.3333333333333333,0.179187,0.537561,0.1044847745542456,0.1044847745542456,3,1,3.574363333333334,10.72309,0.2952944509563076,0.2952944509563076,0.64047944261,5.76431498349,0.03085377413528954,0.03085377413528954)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    for ($c = 1; $c -le $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c - 1]
    }
}
